$wb = $excel.ActiveWorkbook

# --- Settings sheet ---
$ws1 = $wb.Worksheets.Item("Settings")

# Row 2: OrchestratorQueueName / RPAChallenge / description
$ws1.Range("A2").Value = "OrchestratorQueueName"
$ws1.Range("B2").Value = "RPAChallenge"
$ws1.Range("C2").Value = "Orchestrator queue Name. The value must match with the queue name defined on Orchestrator."

# Row 3: OrchestratorQueueFolder / Shared / folder description
$ws1.Range("A3").Value = "OrchestratorQueueFolder"
$ws1.Range("B3").Value = "Shared"
$ws1.Range("C3").Value = "Folder name. The value must match a folder defined in Orchestrator and queue specified as OrchestratorQueueName should be created in this folder. For classic folders leave the value field empty."

# Row 5: ChallengeURL / https://www.rpachallenge.com (description cell cleared)
$ws1.Range("A5").Value = "ChallengeURL"
$ws1.Range("B5").Value = "https://www.rpachallenge.com"
$ws1.Range("C5").Value = ""

# Row 7: ProcessName / chrome (new row of data)
$ws1.Range("A7").Value = "ProcessName"
$ws1.Range("B7").Value = "chrome"

# Remove the 3 trailing empty formatted rows (996-998) so the sheet shrinks to 995 rows
$ws1.Range("A996:A998").EntireRow.Delete()

$ws1.Range("B8").Select()

$wb.Save()
